$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "2024-11-04 10:44:36"
$ws.Range("B14").Value = "High"
$ws.Range("C14").Value = "tasks.py, line 22"
$ws.Range("D14").Value = "Testing!!!"
